$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seed Surveys")

# --- Set R-column values (0/1) for rows with an existing R cell ---
$ws.Range("R2").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("R14").Value = 0
$ws.Range("R15").Value = 0
$ws.Range("R16").Value = 0
$ws.Range("R17").Value = 0
$ws.Range("R18").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("R23").Value = 0
$ws.Range("R25").Value = 0
$ws.Range("R28").Value = 0
$ws.Range("R29").Value = 0
$ws.Range("R31").Value = 0
$ws.Range("R32").Value = 0
$ws.Range("R34").Value = 0
$ws.Range("R35").Value = 0
$ws.Range("R38").Value = 0
$ws.Range("R39").Value = 0
$ws.Range("R41").Value = 0
$ws.Range("R42").Value = 0
$ws.Range("R43").Value = 0
$ws.Range("R45").Value = 0
$ws.Range("R46").Value = 0
$ws.Range("R48").Value = 0
$ws.Range("R51").Value = 0
$ws.Range("R60").Value = 0
$ws.Range("R64").Value = 0
$ws.Range("R72").Value = 1
$ws.Range("R75").Value = 0
$ws.Range("R77").Value = 1
$ws.Range("R86").Value = 1
$ws.Range("R103").Value = 0
$ws.Range("R106").Value = 1
$ws.Range("R108").Value = 0
$ws.Range("R110").Value = 0
$ws.Range("R113").Value = 0
$ws.Range("R114").Value = 0
$ws.Range("R117").Value = 0
$ws.Range("R118").Value = 0
$ws.Range("R121").Value = 0
$ws.Range("R122").Value = 0
$ws.Range("R128").Value = 1
$ws.Range("R130").Value = 0
$ws.Range("R131").Value = 0
$ws.Range("R132").Value = 0
$ws.Range("R133").Value = 0
$ws.Range("R137").Value = 0

# --- Clear R-column values (keep existing style) for rows listed ---
$ws.Range("R66").ClearContents()
$ws.Range("R78").ClearContents()
$ws.Range("R84").ClearContents()
$ws.Range("R85").ClearContents()
$ws.Range("R87").ClearContents()
$ws.Range("R88").ClearContents()
$ws.Range("R123").ClearContents()
$ws.Range("R134").ClearContents()
$ws.Range("R135").ClearContents()
$ws.Range("R136").ClearContents()
$ws.Range("R138").ClearContents()
$ws.Range("R139").ClearContents()
$ws.Range("R140").ClearContents()
$ws.Range("R141").ClearContents()
$ws.Range("R142").ClearContents()
$ws.Range("R143").ClearContents()
$ws.Range("R144").ClearContents()
$ws.Range("R145").ClearContents()
$ws.Range("R148").ClearContents()
$ws.Range("R150").ClearContents()
$ws.Range("R151").ClearContents()

# --- New R43 cell: add value 0 with the centered style used by neighboring R cells ---
$ws.Range("R42").Copy()
$ws.Range("R43").PasteSpecial(-4122)
$ws.Range("R43").Value = 0

# --- Update the view selection to match the saved state ---
$ws.Range("R138").Select()
